$d = $word.ActiveDocument

# --- Paragraph 1: title line (two runs separated by a line break) ---
$d.Content.Find.Execute('🚀המאמר היומי של מייק 28.09.24: ⚡️🚀', $true, $false, $false, $false, $false, $true, 1, $false, '🚀המאמר היומי של מייק 27.09.24: ⚡️🚀', 2) | Out-Null
$d.Content.Find.Execute('Meta-Whisper: Speech-Based Meta-ICL for ASR on Low-Resource Languages', $true, $false, $false, $false, $false, $true, 1, $false, 'REWARD-ROBUST RLHF IN LLMS', 2) | Out-Null

# --- Paragraphs 2-6: body text replacements ---
$d.Content.Find.Execute('מזמן לא סקרתי מאמר על אודיו ומשלים את הפער היום עם סקירה קצרה וקלילה. בדיוק כמו במודלי שפה גם במודלי אודיו כמו whisper למשל יש יכולת למידה in-context או ICL בקצרה. ICL היא יכולת של מודל לבצע משימה שלא אומן עליה באופן מפורש אחרי ש״מראים לו״ כמה דוגמאות המדגימות את המשימה (נגיד, כמה זוגות של שאלות ותשובות רצויות).', $true, $false, $false, $false, $false, $true, 1, $false, 'הסקירה של היום הינה בנושא שהוא די דומה לסקירה של אתמול (26.09.24). נושא של הסקירה הוא שיפור של יישור (alignment) של מודלי שפה במהלך אימון RLHF. גם המאמר הזה מציע שיטה שבאה ״לתקן״ את פונקציית התגמול (reward) אבל מזווית טיפה שונה מאשר המאמר שסקרנו קודם.', 2) | Out-Null
$d.Content.Find.Execute('מתברר שמודלי אודיו גם ניחנים ביכולת כזה. כלומר בהינתן זוג של קטעי אודיו (שאלה ותשובה) ניתן לאמן את המודל לענות על שאלה אחרת, שמוגשת לא לאחר כן בצורה של טקסט. אבל איך ניתן לבחור את הדוגמא מהדאטהסט (אודיו) של שאלות ותשובות שתמקסם את ביצועי המודל לשאלה נתונה.', $true, $false, $false, $false, $false, $true, 1, $false, 'המחברים מצביעים על כך ששימוש בפונקציית תגמול יחידה במהלך אימון RLHF אינו אופטימלי מכמה סיבות. הסיבה הראשונה היא חוסר עקביות בין המתייגים במהלך תיוג הדאטה המשמש לאימון RLHF (כלומר תשובות מועדפות ולא מועדפות לשאלות מהדאטהסט) שעלול לגרום לתשובות ״מבולבלות״ של המודל לאחר האימון. הבעיה השניה היא reward hacking של המודל המתבטא בכך שהמודל לומד להחזיר תשובות הממקסמות את פונקציית התגמול תוך מתן תשובות לא ״מיושרות״ עם העדפות המתייגים או לא הגיוניות.', 2) | Out-Null
$d.Content.Find.Execute('זה בדיוק מה שהמאמר המסוקר עושה. הוא מציע לבחור זוג אודיו (שאלה ותשובה) לשאלה טקסטואלית נתונה על סמך דמיון בין ייצוגה לבין הייצוג של הזוג. הייצוג כאן הוא הפלטים (hidden states) של השכבות השונות של המודל עבור האודיו והשאלה הטקסטואלית. והמטריקה KL divergence הדי סטנדרטי. לדאטהסט אודיו של שאלות ותשובות נתון אני שומרים את כל הפלטים של השכבות ולכל שאלת אודיו בוחרים את הזוג הדומה ביותר לפי מטריקה זו.', $true, $false, $false, $false, $false, $true, 1, $false, 'המאמר ניגש לסוגיה זו מנקודת מבט בייסיאנית. אם נניח שקיימת פונקציית תגמול אידאלית שאין לנו גישה אליה אז ניתן להתבונן בכל פונקציית תגמול שנבנה איזה דגימה ממרחב ״פונקציות תגמול רועשות״. המחברים מציעים לכמת את אי וודאות שיש לנו בפונקציית התגמול על ידי אימון של כמה פונקציות תגמול. ', 2) | Out-Null
$d.Content.Find.Execute('שכחתי לציין שהמודל עובר פיינטיון למשימת ICL בשיטת LoRA הידועה…', $true, $false, $false, $false, $false, $true, 1, $false, 'אז איך כל הסיפור הזה עובד? קודם כל מאמנים פונקציית תגמול רגילה דרך נוסחת Bradley-Terry הסטנדרטי. ', 2) | Out-Null
$d.Content.Find.Execute('זהו זה - סקירה קלילה כמו שהבטחתי.', $true, $false, $false, $false, $false, $true, 1, $false, 'לאחר מכן מאמנים כמה פונקציות תגמול שימדלו לנו את אי הוודאות. בשביל זה לוקחים backbone רגיל (מודל שפה) ומוסיפים אליו כמה ראשים (heads) שכל אחד הוא למעשה פוקנצית תגמול. כל ראש מאומן לפלוט את התוחלת ואת השונות של ערך התגמול והתגמול עצמו מוגרל מהתפלגות גאוסית המוגדרת על ידיהם. ', 2) | Out-Null

# --- Insert two new paragraphs after paragraph 6 (before the URL paragraph) ---
$p6 = $d.Paragraphs(6)
$p6.Range.InsertParagraphAfter()
$p7new = $d.Paragraphs(7)
$p7new.Range.Text = 'פונקציית לוס שהם משתמשים לאימון הראשים היא די לא טריוויאלית אך בגדול ממזערת את השגיאה הריבועית של שערוך התגמול (וזה קצת מורכב ומסתמך על פונקציית תגמול סטנדרטית מהשלב הראשון בנוסף לגישת Bradley Terry). במהלך האימון כל דוגמא מוגרלת (מנווטת) לראש שלו וכך אנו מקבלים כמה פונקציות תגמול. '
$p7new.Range.InsertParagraphAfter()
$p8new = $d.Paragraphs(8)
$p8new.Range.Text = 'המחברים אומרים שהם ״היו רוצים״ (והם השתמשו בה על דוגמאות הצעצוע שלהם) לבנות את הלוס עבור אימון RLHF בתור צירוף לינארי של פונקצית התגמול הרגילה התגמול המינימלי בין כל פונקציות התגמול. כאן האיבר השני למעשה מהווה שערוך של אי הוודאות שדנו בה למעלה. באופן פרקטי במהלך אימון RLHF הם בוחרים ערך התגמול המתקבל בפונקציית התגמול בעלת שונות הנמוכה ביותר. '

# --- Final paragraph: update the URL ---
$d.Content.Find.Execute('https://arxiv.org/abs/2409.10429', $true, $false, $false, $false, $false, $true, 1, $false, 'https://www.arxiv.org/abs/2409.15360', 2) | Out-Null
